# TouchGFX texts.xlsx - "Demo version to be uploaded"
#
# Updates the Translation sheet's toggle-LED related rows (7-10):
#   - Row 7 becomes the "White" toggle button entry
#   - Row 8 becomes the "Orange" toggle button entry
#   - Row 9 becomes the "Toggle LED" entry (previously on row 10)
#   - Row 10 is emptied (its previous content moved up to row 9)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 7: toggleButtonWhite / Label / Center / White Text / LTR
$ws.Range("B7").Value = "toggleButtonWhite"
$ws.Range("C7").Value = "Label"
$ws.Range("D7").Value = "Center"
$ws.Range("E7").Value = "White Text"
$ws.Range("F7").Value = "LTR"

# Row 8: toggleButtonOrange / Label / Center / Orange Text / LTR
$ws.Range("B8").Value = "toggleButtonOrange"
$ws.Range("C8").Value = "Label"
$ws.Range("D8").Value = "Center"
$ws.Range("E8").Value = "Orange Text"
$ws.Range("F8").Value = "LTR"

# Row 9: SingleUseId9 / Label / Center / Toggle LED / LTR
$ws.Range("B9").Value = "SingleUseId9"
$ws.Range("C9").Value = "Label"
$ws.Range("D9").Value = "Center"
$ws.Range("E9").Value = "Toggle LED"
$ws.Range("F9").Value = "LTR"

# Row 10: cleared
$ws.Range("B10:F10").ClearContents()
